$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# Extend formatting/formula pattern by copying the last existing data row (202)
# down into the 4 new rows (203:206), then overwrite the actual data values.
$ws.Range("A202:D202").Copy()
$ws.Range("A203:D206").PasteSpecial()
$excel.CutCopyMode = $false

$newRows = @(
    @{ A = 44105.365277777775; B = 0.36527777777777781; C = 71.099999999999994 },
    @{ A = 44105.364583333336; B = 0.36458333333333331; C = 71.099999999999994 },
    @{ A = 44105.335416666669; B = 0.3354166666666667;  C = 71.099999999999994 },
    @{ A = 44104.918749999997; B = 0.91875000000000007; C = 71.8 }
)

$startRow = 203
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
}

$excel.Goto($ws.Range("A207"), $false)
